# Updated symbol list on Mon Dec 12 03:22:35 UTC 2022 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Hora" (column G) values for the
# cryptocurrency rows (rows 2-51) on the active worksheet. Values are
# written as literal text (matching the source data, which stores them as
# strings) so that formatted numbers such as "283.90" or "0.09200" keep
# their exact textual representation instead of being normalized to a
# floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Price" values (column D) for the rows whose price changed.
$priceUpdates = @{
    2  = "283.90"
    3  = "20.98"
    4  = "6.294"
    5  = "0.06175"
    6  = "3.587"
    7  = "6.560"
    8  = "1.505"
    9  = "0.8201"
    10 = "0.01384"
    11 = "0.1653"
    12 = "0.08462"
    13 = "0.03481"
    14 = "0.03217"
    15 = "0.09200"
    16 = "3.741"
    17 = "0.001668"
    18 = "0.04715"
    19 = "0.006549"
    20 = "0.006182"
    21 = "0.001072"
    22 = "0.0001606"
    23 = "3.840"
    25 = "0.3350"
    26 = "0.1226"
    40 = "0.04732"
    41 = "0.007168"
    42 = "0.004025"
    43 = "0.1103"
    44 = "0.01152"
    45 = "0.00006744"
    46 = "0.00000000753"
    47 = "1.105"
    48 = "0.002850"
    49 = "0.00001908"
    50 = "0.01245"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# The "Hora" column (G) is updated from "2" to "3" for every data row (2-51).
for ($row = 2; $row -le 51; $row++) {
    $cell = $ws.Range("G$row")
    $cell.NumberFormat = "@"
    $cell.Value = "3"
}
